$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (e.g. "54.95" -> 54.95)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.481.00"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "3.404.75"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "562.15"
$ws.Range("E5").Value = "  +3.13%  "
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("E7").Value = "  +2.88%  "
$ws.Range("D8").Value = "3.394.57"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +12.61%  "
$ws.Range("E11").Value = "  +3.47%  "
$ws.Range("D12").Value = "54.95"
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("E13").Value = "  +5.83%  "
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").Value = "3.948.52"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").Value = "3.404.14"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("D19").Value = "65.492.62"
$ws.Range("E19").Value = "  +2.94%  "
$ws.Range("D20").Value = "11.92"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").Value = "472.99"
$ws.Range("E22").Value = "  +14.79%  "
$ws.Range("D23").Value = "5.08"
$ws.Range("E23").Value = "  +17.76%  "
$ws.Range("D24").Value = "4.15"
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").Value = "86.83"
$ws.Range("E25").Value = "  +4.49%  "
$ws.Range("D26").Value = "13.44"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").Value = "10.93"
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("E28").Value = "  +6.63%  "
$ws.Range("D29").Value = "8.88"
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("D30").Value = "31.25"
$ws.Range("E30").Value = "  +7.62%  "
$ws.Range("E31").Value = "  +5.29%  "
$ws.Range("D32").Value = "11.57"
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("D33").Value = "62.62"
$ws.Range("E33").Value = "  +8.17%  "
$ws.Range("D34").Value = "574.32"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -4.44%  "
$ws.Range("E38").Value = "  +3.71%  "
$ws.Range("D39").Value = "35.90"
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("D40").Value = "0.0₃0760"
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("E41").Value = "  +2.12%  "
$ws.Range("D42").Value = "3.092.69"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "2.86"
$ws.Range("E44").Value = "  +2.21%  "
$ws.Range("E45").Value = "  +4.57%  "
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("E47").Value = "  +5.92%  "
$ws.Range("D48").Value = "3.20"
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "137.07"
$ws.Range("E50").Value = "  +4.05%  "
$ws.Range("D51").Value = "8.35"
$ws.Range("E51").Value = "  +3.70%  "
